# Edit viadata.xlsx: update April dates to May ("Mei"), add a 5th data row,
# and highlight Status column cells (green = Pass, red = Fail).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing date cells (F2:F4): April -> Mei (May) ---
$ws.Range("F2").Value = "'28/Mei/2024"
$ws.Range("F3").Value = "'29/Mei/2024"
$ws.Range("F4").Value = "'30/Mei/2024"

# --- Add a 5th row of data (row 5) ---
$ws.Range("A5").Value = "'5"
$ws.Range("B5").Value = "CGK"
$ws.Range("C5").Value = "Jakarta-soekarno Hatta - Indonesia"
$ws.Range("D5").Value = "SIN"
$ws.Range("E5").Value = "Changi Intl Arpt - Singapore"
$ws.Range("F5").NumberFormat = "d-mmm-yy"
$ws.Range("F5").Value = "'23/Mei/2024"
$ws.Range("G5").Value = "AirAsia"
$ws.Range("H5").Value = "Pass"

# give the new row the same thin-border grid as the rest of the table
$ws.Range("A5:H5").Borders.LineStyle = 1

# --- Highlight the Status column: green for Pass, red for Fail ---
$ws.Range("H2").Interior.ColorIndex = 4
$ws.Range("H2").Font.ColorIndex = 1

$ws.Range("H3").Interior.ColorIndex = 4
$ws.Range("H3").Font.ColorIndex = 1

$ws.Range("H4").Interior.ColorIndex = 3
$ws.Range("H4").Font.ColorIndex = 1

$ws.Range("H5").Interior.ColorIndex = 4
$ws.Range("H5").Font.ColorIndex = 1

# --- Update the selected cell shown in the sheet view ---
[void]$ws.Range("F6").Select()
